$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: RandomForestRegressor - name unchanged, update values
$ws.Range("B3").Value = 64990492163129.3
$ws.Range("C3").Value = 63815704116543.59
$ws.Range("D3").Value = 19694966908024.94

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 20369970566000.13
$ws.Range("C4").Value = 20369970566000.13
$ws.Range("D4").Value = 18481125231428.3

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 120646324675591.2
$ws.Range("C5").Value = 120373777549764.1
$ws.Range("D5").Value = 128381669412927.2
